$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "model_26_9_24"
$ws.Range("B2").Value = 0.9975182599622946
$ws.Range("C2").Value = 0.6716090796131864
$ws.Range("D2").Value = 0.9952823018316802
$ws.Range("E2").Value = 0.9815534291385378
$ws.Range("F2").Value = 0.992625800406205
$ws.Range("G2").Value = 0.00100038697235429
$ws.Range("H2").Value = 0.194947036090733
$ws.Range("I2").Value = 0.002952246237628282
$ws.Range("J2").Value = 0.003675796468756666
$ws.Range("K2").Value = 0.003314026252657492
$ws.Range("L2").Value = 0.004709588608930119
$ws.Range("M2").Value = 0.0316288945800243
$ws.Range("N2").Value = 1.000677330796317
$ws.Range("O2").Value = 0.03188294600726131
$ws.Range("P2").Value = 1179.814736762965
$ws.Range("Q2").Value = 2833.367088515676
$ws.Range("A3").Value = "model_26_9_23"
$ws.Range("B3").Value = 0.9975380852147034
$ws.Range("C3").Value = 0.6715086549251414
$ws.Range("D3").Value = 0.9953258533788878
$ws.Range("E3").Value = 0.9816897684014174
$ws.Range("F3").Value = 0.9926863488841486
$ws.Range("G3").Value = 0.0009923954325748708
$ws.Range("H3").Value = 0.1950066525236772
$ws.Range("I3").Value = 0.002924992503540374
$ws.Range("J3").Value = 0.003648628526009505
$ws.Range("K3").Value = 0.003286815266175269
$ws.Range("L3").Value = 0.004712475284367118
$ws.Range("M3").Value = 0.03150230836898894
$ws.Range("N3").Value = 1.000671919974153
$ws.Range("O3").Value = 0.03175534302317685
$ws.Range("P3").Value = 1179.830777817121
$ws.Range("Q3").Value = 2833.383129569833
$ws.Range("A4").Value = "model_26_9_22"
$ws.Range("B4").Value = 0.9975597949208335
$ws.Range("C4").Value = 0.6713977041321872
$ws.Range("D4").Value = 0.9953727613465501
$ws.Range("E4").Value = 0.9818424356344664
$ws.Range("F4").Value = 0.9927528708599289
$ws.Range("G4").Value = 0.000983644271350873
$ws.Range("H4").Value = 0.1950725177071994
$ws.Range("I4").Value = 0.002895638385047672
$ws.Range("J4").Value = 0.00361820695441493
$ws.Range("K4").Value = 0.003256919740388325
$ws.Range("L4").Value = 0.004715788383673408
$ws.Range("M4").Value = 0.03136310366259808
$ws.Range("N4").Value = 1.000665994836017
$ws.Range("O4").Value = 0.03161502018873232
$ws.Range("P4").Value = 1179.848492478257
$ws.Range("Q4").Value = 2833.400844230968
$ws.Range("A5").Value = "model_26_9_21"
$ws.Range("B5").Value = 0.9975847723161132
$ws.Range("C5").Value = 0.671267092698471
$ws.Range("D5").Value = 0.9954258006110998
$ws.Range("E5").Value = 0.9820219051179624
$ws.Range("F5").Value = 0.9928295709913497
$ws.Range("G5").Value = 0.0009735759078391489
$ws.Range("H5").Value = 0.1951500542963736
$ws.Range("I5").Value = 0.00286244741698936
$ws.Range("J5").Value = 0.003582444573502048
$ws.Range("K5").Value = 0.003222450067323766
$ws.Range("L5").Value = 0.004718747327970767
$ws.Range("M5").Value = 0.03120217793422678
$ws.Range("N5").Value = 1.000659177861323
$ws.Range("O5").Value = 0.03145280186346458
$ws.Range("P5").Value = 1179.869069524075
$ws.Range("Q5").Value = 2833.421421276787
$ws.Range("A6").Value = "model_26_9_20"
$ws.Range("B6").Value = 0.9976120826822439
$ws.Range("C6").Value = 0.67112375251888
$ws.Range("D6").Value = 0.9954846426636684
$ws.Range("E6").Value = 0.9822169256623632
$ws.Range("F6").Value = 0.9929137674100295
$ws.Range("G6").Value = 0.0009625671260681696
$ws.Range("H6").Value = 0.1952351472189521
$ws.Range("I6").Value = 0.002825625173998755
$ws.Range("J6").Value = 0.003543583376273183
$ws.Range("K6").Value = 0.003184611500800672
$ws.Range("L6").Value = 0.004721933152998506
$ws.Range("M6").Value = 0.03102526593065996
$ws.Range("N6").Value = 1.000651724158776
$ws.Range("O6").Value = 0.03127446885712171
$ws.Range("P6").Value = 1179.891813505729
$ws.Range("Q6").Value = 2833.444165258441
$ws.Range("A7").Value = "model_26_9_19"
$ws.Range("B7").Value = 0.9976424503088017
$ws.Range("C7").Value = 0.6709585797988461
$ws.Range("D7").Value = 0.9955499156627761
$ws.Range("E7").Value = 0.982435808269581
$ws.Range("F7").Value = 0.9930077599437028
$ws.Range("G7").Value = 0.0009503259656209959
$ws.Range("H7").Value = 0.1953332008806542
$ws.Range("I7").Value = 0.002784778566361065
$ws.Range("J7").Value = 0.003499967252673531
$ws.Range("K7").Value = 0.003142370479224665
$ws.Range("L7").Value = 0.004725410565675619
$ws.Range("M7").Value = 0.03082735742195552
$ws.Range("N7").Value = 1.000643436051091
$ws.Range("O7").Value = 0.03107497069630434
$ws.Range("P7").Value = 1179.917411021032
$ws.Range("Q7").Value = 2833.469762773744
$ws.Range("A8").Value = "model_26_9_18"
$ws.Range("B8").Value = 0.9976772165361729
$ws.Range("C8").Value = 0.6707734207018483
$ws.Range("D8").Value = 0.9956236854125367
$ws.Range("E8").Value = 0.9826910200937821
$ws.Range("F8").Value = 0.9931156913971628
$ws.Range("G8").Value = 0.0009363117335049415
$ws.Range("H8").Value = 0.1954431193190886
$ws.Range("I8").Value = 0.002738614852954408
$ws.Range("J8").Value = 0.003449111907838505
$ws.Range("K8").Value = 0.003093865191877297
$ws.Range("L8").Value = 0.00472812256597563
$ws.Range("M8").Value = 0.03059921132161647
$ws.Range("N8").Value = 1.000633947451918
$ws.Range("O8").Value = 0.03084499206772871
$ws.Range("P8").Value = 1179.947124176637
$ws.Range("Q8").Value = 2833.499475929349
$ws.Range("A9").Value = "model_26_9_17"
$ws.Range("B9").Value = 0.9977153222114611
$ws.Range("C9").Value = 0.6705623074267818
$ws.Range("D9").Value = 0.9957045985630982
$ws.Range("E9").Value = 0.9829733105684954
$ws.Range("F9").Value = 0.9932346000026666
$ws.Range("G9").Value = 0.0009209513732125989
$ws.Range("H9").Value = 0.1955684452787868
$ws.Range("I9").Value = 0.002687980934505769
$ws.Range("J9").Value = 0.003392860676218941
$ws.Range("K9").Value = 0.003040426681664194
$ws.Range("L9").Value = 0.004730785327554941
$ws.Range("M9").Value = 0.03034718064685085
$ws.Range("N9").Value = 1.00062354743137
$ws.Range("O9").Value = 0.03059093701767329
$ws.Range("P9").Value = 1179.980206641835
$ws.Range("Q9").Value = 2833.532558394547
$ws.Range("A10").Value = "model_26_9_16"
$ws.Range("B10").Value = 0.9977595861091405
$ws.Range("C10").Value = 0.6703297350277386
$ws.Range("D10").Value = 0.9957957754528541
$ws.Range("E10").Value = 0.9833116777892161
$ws.Range("F10").Value = 0.9933731260102349
$ws.Range("G10").Value = 0.0009031086395211981
$ws.Range("H10").Value = 0.1957065103014633
$ws.Range("I10").Value = 0.002630924162296774
$ws.Range("J10").Value = 0.003325435188608842
$ws.Range("K10").Value = 0.00297817194880567
$ws.Range("L10").Value = 0.004732127632102942
$ws.Range("M10").Value = 0.03005176599671304
$ws.Range("N10").Value = 1.00061146667327
$ws.Range("O10").Value = 0.030293149520982
$ws.Range("P10").Value = 1180.019335404473
$ws.Range("Q10").Value = 2833.571687157185
$ws.Range("A11").Value = "model_26_9_15"
$ws.Range("B11").Value = 0.9978062115700133
$ws.Range("C11").Value = 0.6700643227973577
$ws.Range("D11").Value = 0.9958919302905807
$ws.Range("E11").Value = 0.9836735635044036
$ws.Range("F11").Value = 0.9935202792060389
$ws.Range("G11").Value = 0.0008843139620253907
$ws.Range("H11").Value = 0.1958640704666288
$ws.Range("I11").Value = 0.002570752284448824
$ws.Range("J11").Value = 0.003253323236530066
$ws.Range("K11").Value = 0.002912040086241603
$ws.Range("L11").Value = 0.004734289055704243
$ws.Range("M11").Value = 0.02973741686874283
$ws.Range("N11").Value = 1.00059874138373
$ws.Range("O11").Value = 0.02997627545985583
$ws.Range("P11").Value = 1180.061396795438
$ws.Range("Q11").Value = 2833.61374854815
$ws.Range("A12").Value = "model_26_9_14"
$ws.Range("B12").Value = 0.997859006336132
$ws.Range("C12").Value = 0.6697643456087748
$ws.Range("D12").Value = 0.9959997195865378
$ws.Range("E12").Value = 0.9840903579374524
$ws.Range("F12").Value = 0.993687728384444
$ws.Range("G12").Value = 0.000863032443642609
$ws.Range("H12").Value = 0.1960421498841114
$ws.Range("I12").Value = 0.002503299782806605
$ws.Range("J12").Value = 0.003170269778155372
$ws.Range("K12").Value = 0.002836787041329812
$ws.Range("L12").Value = 0.004734775369115369
$ws.Range("M12").Value = 0.02937741383516611
$ws.Range("N12").Value = 1.000584332331842
$ws.Range("O12").Value = 0.02961338078919523
$ws.Range("P12").Value = 1180.110116547127
$ws.Range("Q12").Value = 2833.662468299839
$ws.Range("A13").Value = "model_26_9_13"
$ws.Range("B13").Value = 0.9979159433091874
$ws.Range("C13").Value = 0.669419041080241
$ws.Range("D13").Value = 0.996117412979891
$ws.Range("E13").Value = 0.9845425867811534
$ws.Range("F13").Value = 0.993869918677271
$ws.Range("G13").Value = 0.0008400812057109638
$ws.Range("H13").Value = 0.196247137568628
$ws.Range("I13").Value = 0.002429649484435632
$ws.Range("J13").Value = 0.00308015540409473
$ws.Range("K13").Value = 0.002754909217746575
$ws.Range("L13").Value = 0.004734816247618888
$ws.Range("M13").Value = 0.02898415439013123
$ws.Range("N13").Value = 1.00056879276496
$ws.Range("O13").Value = 0.02921696258301451
$ws.Range("P13").Value = 1180.164023994668
$ws.Range("Q13").Value = 2833.71637574738
$ws.Range("A14").Value = "model_26_9_12"
$ws.Range("B14").Value = 0.997982505668534
$ws.Range("C14").Value = 0.6690360496881356
$ws.Range("D14").Value = 0.9962511106226944
$ws.Range("E14").Value = 0.9850867022146867
$ws.Range("F14").Value = 0.9940836682807568
$ws.Range("G14").Value = 0.0008132499840165691
$ws.Range("H14").Value = 0.1964744978033484
$ws.Range("I14").Value = 0.002345984029617712
$ws.Range("J14").Value = 0.002971731046841627
$ws.Range("K14").Value = 0.002658848379083785
$ws.Range("L14").Value = 0.004730843814423458
$ws.Range("M14").Value = 0.02851753818295978
$ws.Range("N14").Value = 1.000550626182169
$ws.Range("O14").Value = 0.02874659839429081
$ws.Range("P14").Value = 1180.228944024529
$ws.Range("Q14").Value = 2833.781295777241
$ws.Range("A15").Value = "model_26_9_11"
$ws.Range("B15").Value = 0.9980584223544827
$ws.Range("C15").Value = 0.6685721883210836
$ws.Range("D15").Value = 0.9964003230745442
$ws.Range("E15").Value = 0.9857244076393253
$ws.Range("F15").Value = 0.9943288991926137
$ws.Range("G15").Value = 0.0007826480424539799
$ws.Range("H15").Value = 0.1967498659485983
$ws.Range("I15").Value = 0.002252609700895514
$ws.Range("J15").Value = 0.002844657274399132
$ws.Range("K15").Value = 0.002548639580214189
$ws.Range("L15").Value = 0.00472205732615566
$ws.Range("M15").Value = 0.02797584748410636
$ws.Range("N15").Value = 1.000529906562641
$ws.Range("O15").Value = 0.02820055669623263
$ws.Range("P15").Value = 1180.30565492362
$ws.Range("Q15").Value = 2833.858006676332
$ws.Range("A16").Value = "model_26_9_10"
$ws.Range("B16").Value = 0.9981412037852979
$ws.Range("C16").Value = 0.6680417155851922
$ws.Range("D16").Value = 0.9965639484545433
$ws.Range("E16").Value = 0.9864286791316127
$ws.Range("F16").Value = 0.9945989586232206
$ws.Range("G16").Value = 0.0007492789289762485
$ws.Range("H16").Value = 0.1970647774798524
$ws.Range("I16").Value = 0.002150216034482764
$ws.Range("J16").Value = 0.002704319068244827
$ws.Range("K16").Value = 0.002427272639785521
$ws.Range("L16").Value = 0.004708527868121364
$ws.Range("M16").Value = 0.02737295981395232
$ws.Range("N16").Value = 1.000507313377375
$ws.Range("O16").Value = 0.02759282647703918
$ws.Range("P16").Value = 1180.392798483868
$ws.Range("Q16").Value = 2833.94515023658
$ws.Range("A17").Value = "model_26_9_9"
$ws.Range("B17").Value = 0.9982313656914039
$ws.Range("C17").Value = 0.6674287491262394
$ws.Range("D17").Value = 0.9967381400204455
$ws.Range("E17").Value = 0.9872228081053046
$ws.Range("F17").Value = 0.9948963041037656
$ws.Range("G17").Value = 0.0007129347531557891
$ws.Range("H17").Value = 0.1974286608486591
$ws.Range("I17").Value = 0.002041210248882864
$ws.Range("J17").Value = 0.002546075213646777
$ws.Range("K17").Value = 0.002293643122967967
$ws.Range("L17").Value = 0.004695535028056486
$ws.Range("M17").Value = 0.0267008380609259
$ws.Range("N17").Value = 1.00048270587025
$ws.Range("O17").Value = 0.02691530606898859
$ws.Range("P17").Value = 1180.492241304075
$ws.Range("Q17").Value = 2834.044593056787
$ws.Range("A18").Value = "model_26_9_8"
$ws.Range("B18").Value = 0.998330734434299
$ws.Range("C18").Value = 0.6667061830821097
$ws.Range("D18").Value = 0.9969265991643285
$ws.Range("E18").Value = 0.9881256664747323
$ws.Range("F18").Value = 0.9952276852544737
$ws.Range("G18").Value = 0.0006728793104659488
$ws.Range("H18").Value = 0.1978576072656823
$ws.Range("I18").Value = 0.001923276082977222
$ws.Range("J18").Value = 0.002366165157135224
$ws.Range("K18").Value = 0.002144717694639885
$ws.Range("L18").Value = 0.004676633383263123
$ws.Range("M18").Value = 0.02593991731802453
$ws.Range("N18").Value = 1.000455585580159
$ws.Range("O18").Value = 0.02614827341470634
$ws.Range("P18").Value = 1180.607889150122
$ws.Range("Q18").Value = 2834.160240902834
$ws.Range("A19").Value = "model_26_9_7"
$ws.Range("B19").Value = 0.9984451527691709
$ws.Range("C19").Value = 0.6658769603814604
$ws.Range("D19").Value = 0.9971382025599234
$ws.Range("E19").Value = 0.9891935514566719
$ws.Range("F19").Value = 0.9956117419667141
$ws.Range("G19").Value = 0.0006267573920275544
$ws.Range("H19").Value = 0.1983498696813433
$ws.Range("I19").Value = 0.00179085868232419
$ws.Range("J19").Value = 0.002153370710127606
$ws.Range("K19").Value = 0.001972119433542498
$ws.Range("L19").Value = 0.004643888758097825
$ws.Range("M19").Value = 0.02503512316781275
$ws.Range("N19").Value = 1.000424357868676
$ws.Range("O19").Value = 0.02523621172485179
$ws.Range("P19").Value = 1180.74990205355
$ws.Range("Q19").Value = 2834.302253806262
$ws.Range("A20").Value = "model_26_9_6"
$ws.Range("B20").Value = 0.9985712455960988
$ws.Range("C20").Value = 0.6648915899405275
$ws.Range("D20").Value = 0.9973690516320017
$ws.Range("E20").Value = 0.9904007114271541
$ws.Range("F20").Value = 0.9960401093065364
$ws.Range("G20").Value = 0.0005759294973047023
$ws.Range("H20").Value = 0.1989348281408677
$ws.Range("I20").Value = 0.001646397701526423
$ws.Range("J20").Value = 0.001912823326548944
$ws.Range("K20").Value = 0.001779607610137752
$ws.Range("L20").Value = 0.004594761424558987
$ws.Range("M20").Value = 0.02399853114889956
$ws.Range("N20").Value = 1.000389943887528
$ws.Range("O20").Value = 0.02419129353187011
$ws.Range("P20").Value = 1180.919050610545
$ws.Range("Q20").Value = 2834.471402363257
$ws.Range("A21").Value = "model_26_9_5"
$ws.Range("B21").Value = 0.9987030318834146
$ws.Range("C21").Value = 0.6637321856160763
$ws.Range("D21").Value = 0.9976135116026333
$ws.Range("E21").Value = 0.9916941794242169
$ws.Range("F21").Value = 0.9964970610966682
$ws.Range("G21").Value = 0.0005228065742899436
$ws.Range("H21").Value = 0.1996231006315209
$ws.Range("I21").Value = 0.001493419277982022
$ws.Range("J21").Value = 0.001655077584439994
$ws.Range("K21").Value = 0.001574249698484545
$ws.Range("L21").Value = 0.004528532791504232
$ws.Range("M21").Value = 0.02286496390309732
$ws.Range("N21").Value = 1.000353976014352
$ws.Range("O21").Value = 0.0230486211820012
$ws.Range("P21").Value = 1181.112598002281
$ws.Range("Q21").Value = 2834.664949754992
$ws.Range("A22").Value = "model_26_9_4"
$ws.Range("B22").Value = 0.998835058325588
$ws.Range("C22").Value = 0.6623626221371663
$ws.Range("D22").Value = 0.9978477854560971
$ws.Range("E22").Value = 0.9930856100188789
$ws.Range("F22").Value = 0.9969686679419353
$ws.Range("G22").Value = 0.0004695868450878726
$ws.Range("H22").Value = 0.2004361326746636
$ws.Range("I22").Value = 0.001346815133802643
$ws.Range("J22").Value = 0.001377811110102284
$ws.Range("K22").Value = 0.001362305683914716
$ws.Range("L22").Value = 0.004438765344917176
$ws.Range("M22").Value = 0.02166995258619346
$ws.Range("N22").Value = 1.000317942596728
$ws.Range("O22").Value = 0.02184401122642675
$ws.Range("P22").Value = 1181.327314605717
$ws.Range("Q22").Value = 2834.879666358428
$ws.Range("A23").Value = "model_26_9_3"
$ws.Range("B23").Value = 0.9989689466018259
$ws.Range("C23").Value = 0.6607327506967184
$ws.Range("D23").Value = 0.9980874386808062
$ws.Range("E23").Value = 0.9945607585762859
$ws.Range("F23").Value = 0.9974625380685743
$ws.Range("G23").Value = 0.000415616612402575
$ws.Range("H23").Value = 0.2014036947684938
$ws.Range("I23").Value = 0.001196844680895394
$ws.Range("J23").Value = 0.00108386239199466
$ws.Range("K23").Value = 0.001140356366667842
$ws.Range("L23").Value = 0.004312282202522755
$ws.Range("M23").Value = 0.02038667732619946
$ws.Range("N23").Value = 1.00028140103662
$ws.Range("O23").Value = 0.02055042836903902
$ws.Range("P23").Value = 1181.571494654887
$ws.Range("Q23").Value = 2835.123846407598
$ws.Range("A24").Value = "model_26_9_2"
$ws.Range("B24").Value = 0.9990951044454108
$ws.Range("C24").Value = 0.6588596541671699
$ws.Range("D24").Value = 0.9983211425179247
$ws.Range("E24").Value = 0.9960420823072649
$ws.Range("F24").Value = 0.9979536669418103
$ws.Range("G24").Value = 0.0003647625095291208
$ws.Range("H24").Value = 0.2025156457819907
$ws.Range("I24").Value = 0.001050597242158066
$ws.Range("J24").Value = 0.0007886831643587288
$ws.Range("K24").Value = 0.0009196389913595135
$ws.Range("L24").Value = 0.004154183685129559
$ws.Range("M24").Value = 0.01909875675349369
$ws.Range("N24").Value = 1.00024696931075
$ws.Range("O24").Value = 0.01925216288659143
$ws.Range("P24").Value = 1181.832528149953
$ws.Range("Q24").Value = 2835.384879902665
$ws.Range("A25").Value = "model_26_9_1"
$ws.Range("B25").Value = 0.9992028092581945
$ws.Range("C25").Value = 0.6567678479050425
$ws.Range("D25").Value = 0.9985403114811828
$ws.Range("E25").Value = 0.9974358060318398
$ws.Range("F25").Value = 0.9984152448504306
$ws.Range("G25").Value = 0.0003213468052524113
$ws.Range("H25").Value = 0.2037574323405151
$ws.Range("I25").Value = 0.0009134454524296129
$ws.Range("J25").Value = 0.0005109597444510161
$ws.Range("K25").Value = 0.0007122020638180409
$ws.Range("L25").Value = 0.003933210779947714
$ws.Range("M25").Value = 0.01792614864527267
$ws.Range("N25").Value = 1.000217573892414
$ws.Range("O25").Value = 0.01807013608804188
$ws.Range("P25").Value = 1182.085979255683
$ws.Range("Q25").Value = 2835.638331008394
$ws.Range("A26").Value = "model_26_9_0"
$ws.Range("B26").Value = 0.9992672939283358
$ws.Range("C26").Value = 0.6544395980212336
$ws.Range("D26").Value = 0.9987042288268448
$ws.Range("E26").Value = 0.998582808586285
$ws.Range("F26").Value = 0.9987836603809076
$ws.Range("G26").Value = 0.0002953530980366107
$ws.Range("H26").Value = 0.2051395820467021
$ws.Range("I26").Value = 0.0008108690794301227
$ws.Range("J26").Value = 0.0002823997605413487
$ws.Range("K26").Value = 0.0005466330790953092
$ws.Range("L26").Value = 0.003668562820852012
$ws.Range("M26").Value = 0.01718584004454279
$ws.Range("N26").Value = 1.000199974364537
$ws.Range("O26").Value = 0.01732388113796566
$ws.Range("P26").Value = 1182.254677949784
$ws.Range("Q26").Value = 2835.807029702496
